$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.714.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.293.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.32%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "268.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.80%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0933"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.38%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.642.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.845"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.319.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.769.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000105"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.96%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.66%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0884"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0348"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.235"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.17%  "
$ws.Range("E46").Value = "  +5.59%  "
$ws.Range("E47").Value = "  +3.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "99.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.521.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.428"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.45%  "
